$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add time-varying selectivity and catchability (q) parameters to the
# build params / map / data metadata sheet.

$ws.Range("F19").Value = "Time_varying_sel"
$ws.Range("G19").Value = "Wether a time-varying selectivity should be estimated for logistic, double logistic selectivity. 0 = no, 1 = random walk from mean selectivity following Dorn 2018, 2 = random effect."

$ws.Range("F20").Value = "Sel_sd_prior"
$ws.Range("G20").Value = "The sd to use for the random walk of time varying selectivity if set to 1"

$ws.Range("F21").Value = "Time_varying_q"

$ws.Range("F22").Value = "Q_sd_prior"
$ws.Range("G22").Value = "The sd to use for the random walk of time varying q if set to 1"

$ws.Range("G21").Value = "Wether a time-varying q should be estimated. 0 = no, 1 = random walk from mean selectivity following Dorn 2018, 2 = random effect."

$ws.Range("F23").Value = "Selectivity_index"
$ws.Range("G23").Value = "index to use if selectivitys of different surveys are to be the same"

$ws.Range("F24").Value = "Q_index"
$ws.Range("G24").Value = "index to use if catchability coefficients are to be set the same"

# Keep the view consistent with the saved workbook (top-left cell and
# the active selection after the edits were made).
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("F20").Select() | Out-Null
